$wb = $excel.ActiveWorkbook

# --- 1. Swap the "LD" and "GA3" sheet names -------------------------------
# Worksheets are addressed positionally (index 1 and 3) since the rename
# itself changes the .Name we would otherwise look up by.
$wsOne   = $wb.Worksheets.Item(1)   # currently "LD"
$wsThree = $wb.Worksheets.Item(3)   # currently "GA3"

# Use a temporary name to dodge the duplicate-name collision while swapping.
$wsOne.Name   = "__tmp_swap__"
$wsThree.Name = "LD"
$wsOne.Name   = "GA3"

# --- 2. Swap the "Proportion" (column E) data between the two sheets -----
# sheet1 (now "GA3") held the old "LD" sheet's proportions and vice-versa;
# the data itself needs to move along with (opposite to) the name swap.
$rangeOne   = $wsOne.Range("E2:E316")
$rangeThree = $wsThree.Range("E2:E316")

$valuesOne   = $rangeOne.Value2
$valuesThree = $rangeThree.Value2

$rangeOne.Value2   = $valuesThree
$rangeThree.Value2 = $valuesOne

# --- 3. Rename a column header and the three proportion labels -----------
# These text values are shared across all three worksheets, so update each
# sheet individually.
foreach ($ws in $wb.Worksheets) {
    $ws.Range("D1").Value = "Segment"
    $ws.Range("D2:D106").Value = "Prop_1"
    $ws.Range("D107:D211").Value = "Prop_2"
    $ws.Range("D212:D316").Value = "Prop_3"
}
